$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B61").Value = 'Aiden Patel'
$ws.Range("C61").Value = '11/12/2024'
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 40
$ws.Range("F61").Value = 'created level1_pvp scene and created some duplicates of scripts with slight modifications to get the bones of PvP in'

$ws.Range("B62").Value = 'Aiden Patel'
$ws.Range("C62").Value = '11/13/2024'
$ws.Range("D62").Value = 2
$ws.Range("E62").Value = 5
$ws.Range("F62").Value = 'Added controller support, controls could be modified to be more intuitive. Added every PvP stage and connected it to the game with menus so that it''s playable in theory. Mines still don''t work since they aren''t designed to collide with Player tanks by default so that still needs to be modified.'

$ws.Range("B63").Value = 'John Newman'
$ws.Range("C63").Value = '11/13/2024'
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 19
$ws.Range("F63").Value = 'Warning fixes Changed a few lines to deal with warnings, mainly unused deltas. enemy_tank.gd: changed code so mine isn''t shadowing the function call hole.tscn: saved it or something so there wasn''t a warning about loading hole.png instead of uid whatever'

$ws.Range("B64").Value = 'John Newman'
$ws.Range("C64").Value = '11/19/2024'
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = 'Added new folder for audio files, added demo music'

$ws.Range("B65").Value = 'Aiden Patel'
$ws.Range("C65").Value = '11/20/2024'
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 30
$ws.Range("F65").Value = 'improve controller controls. Now steer tank left and right with left thumb stick and move forward and back with right and left triggers. Shoot with right bumper and bottom face button (X) and mine with left bumper and right face button (O). Improved aiming with controller to be global based rotation so that its logical that moving the thumb stick right points to the right of the screen, etc'

$ws.Range("B66").Value = 'Aiden Patel'
$ws.Range("C66").Value = '11/20/2024'
$ws.Range("D66").Value = 0
$ws.Range("E66").Value = 40
$ws.Range("F66").Value = 'Added menu navigation for keyboard and controller so that UI buttons are highlighted. Can’t actually select anything on controller rn but you can move it around at least'

$ws.Range("B67").Value = 'John Newman'
$ws.Range("C67").Value = '11/20/2024'
$ws.Range("D67").Value = 0
$ws.Range("E67").Value = 2
$ws.Range("F67").Value = 'Added the game music to audio folder'

$ws.Range("B68").Value = 'Andrew McFerrin'
$ws.Range("C68").Value = '11/22/2024'
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = 50
$ws.Range("F68").Value = 'Added basic tank treads and started to create objects to decorate background'

$ws.Range("B69").Value = 'Andrew McFerrin'
$ws.Range("C69").Value = '11/22/2024'
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 40
$ws.Range("F69").Value = 'created some more background details for map'

$ws.Range("B70").Value = 'Aiden Patel'
$ws.Range("C70").Value = '11/23/2024'
$ws.Range("D70").Value = 1
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 'added audio buses and made sliders in the option menu usable'

$ws.Range("B71").Value = 'Landon Pyko'
$ws.Range("C71").Value = '11/23/2024'
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 'Added a couple audio files for explosion and tank destruction. Tank destruction is not playing properly for some reason. Added an audio manager so that the game music plays between scenes and doesn''t reset between scenes. Created the Final Sprint Requirements Artifacts. Feel free to add more if you think it is necessary'

$ws.Range("B72").Value = 'Kai Achen'
$ws.Range("C72").Value = '11/23/2024'
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 30
$ws.Range("F72").Value = 'changed bullets to only ricochet twice before being destroyed (amount can be changed by changing the ricochet_bank value in the bullet.tscn script'

$ws.Range("B73").Value = 'Kai Achen'
$ws.Range("C73").Value = '11/23/2024'
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 'added sounds for placing mines, player shoot, enemy shoot, bullet ricochet'

$ws.Range("B74").Value = 'Kai Achen'
$ws.Range("C74").Value = '11/23/2024'
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = 'discovered bug bullets will destroy broken walls only after ricocheting off another wall'

$ws.Range("B75").Value = 'Kai Achen'
$ws.Range("C75").Value = '11/23/2024'
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = "x"
$ws.Range("F75").Font.Italic = $true
$ws.Range("F75").Font.Name = "Aptos Narrow"
$ws.Range("F75").Font.Size = 11
$ws.Range("F75").Font.Italic = $false

$ws.Range("F75").Value = 'discovered bug player 1 can only shoot five bullets total in PVP mode (unable to test player 2 at the moment but presumed the same)'
$ws.Range("F75").Characters(53, 5).Font.Italic = $true
$ws.Range("F75").Characters(58, 74).Font.Italic = $false

$ws.Range("B76").Value = 'Kai Achen'
$ws.Range("C76").Value = '11/23/2024'
$ws.Range("D76").Value = 0
$ws.Range("E76").Value = 10
$ws.Range("F76").Value = 'centered game mode select menu'

$ws.Range("B77").Value = 'Kai Achen'
$ws.Range("C77").Value = '11/23/2024'
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 'fixed volume sliders to work and adjusted music and SFX buses to correspond to correct sliders'

$ws.Range("B78").Value = 'Aiden Patel'
$ws.Range("C78").Value = '11/24/2024'
$ws.Range("D78").Value = 2
$ws.Range("E78").Value = 0
$ws.Range("F78").Value = 'refined ricochet bank and fixed bugs, added music speed up with each level, added mine collision with players, added player and player collision so players can’t drive through each other'

$ws.Range("B79").Value = 'John Newman'
$ws.Range("C79").Value = '11/24/2024'
$ws.Range("D79").Value = 0
$ws.Range("E79").Value = 51
$ws.Range("F79").Value = 'checked for errors and suggestions to the game process.'

$ws.Range("B80").Value = 'John Newman'
$ws.Range("C80").Value = '11/24/2024'
$ws.Range("D80").Value = 0
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = 'checked what requirements are left'

$ws.Range("B81").Value = 'John Newman'
$ws.Range("C81").Value = '11/24/2024'
$ws.Range("D81").Value = 0
$ws.Range("E81").Value = 45
$ws.Range("F81").Value = 'Added PVP ready screen Added pvp ready screen. Players mark ready and then wait 3 seconds, in which they can back out. Other details can be added later'

$ws.Range("B82").Value = 'Aiden Patel'
$ws.Range("C82").Value = '11/24/2024'
$ws.Range("D82").Value = 1
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 'added controller selecting for ready up PvP screen, added sounds for mines destroying players and walls'

$ws.Range("B83").Value = 'Aiden Patel'
$ws.Range("C83").Value = '11/24/2024'
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 25
$ws.Range("F83").Value = 'updated worklog'

$excel.ActiveWindow.ScrollRow = 68
$ws.Range("E84").Select()
